$wb = $excel.ActiveWorkbook

# The QREST "Level 1 Data Review" SOP workbook is being repurposed to document
# the Level 2 Data Review instead (the Level 1 SOP is being edited separately
# to cover only checks internal to the monitoring network, while this sheet
# now documents the Level 2 review of data from outside/downwind sites).
# Rename the main worksheet tab accordingly.
$ws = $wb.Worksheets.Item("Level1Review")
$ws.Name = "Level2Review"

# Restore a sensible cursor/selection position on that sheet (B18), matching
# where the author's review left off after editing.
$ws.Activate() | Out-Null
$ws.Range("B18").Select() | Out-Null
